$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header "Injected Bugs Stats" (row 6) - same underline style as "Codebase Stats" (A1)
$ws.Range("A6").Value = "Injected Bugs Stats"
$ws.Range("A6").Font.Underline = $true

# Summary rows
$ws.Range("A7").Value = "Number of unique bug id's"
$ws.Range("B7").Value = 69

$ws.Range("A8").Value = "Total number of alerts flagged in comments"
$ws.Range("B8").Value = 372

$ws.Range("A9").Value = "Number of alerts associated with each CERT Code:"

# Per-CERT-code counts
$ws.Range("A10").Value = "MEM35-C"
$ws.Range("B10").Value = 52

$ws.Range("A11").Value = "ARR30-C"
$ws.Range("B11").Value = 22

$ws.Range("A12").Value = "INT31-C"
$ws.Range("B12").Value = 1

# Note: EXP34-C's shared string is registered before EXP33-C's (matching the
# original authoring order), so write A14 first, then A13.
$ws.Range("A14").Value = "EXP34-C"
$ws.Range("B14").Value = 15

$ws.Range("A13").Value = "EXP33-C"
$ws.Range("B13").Value = 17

$ws.Range("A15").Value = "MSC21-C"
$ws.Range("B15").Value = 32

$ws.Range("A16").Value = "MEM10-C"
$ws.Range("B16").Value = 7

$ws.Range("A17").Value = "EXP08-C"
$ws.Range("B17").Value = 3

$ws.Range("A18").Value = "MEM00-C"
$ws.Range("B18").Value = 2

$ws.Range("A19").Value = "MEM01-C"
$ws.Range("B19").Value = 2

$ws.Range("A20").Value = "MEM30-C"
$ws.Range("B20").Value = 2

# Move selection to A21, matching the post-edit cursor position
$ws.Range("A21").Select() | Out-Null
